$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as Text so values like "4.00" or "1.01"
# are not silently coerced into numbers by Excel (matches source data which
# stores prices/volumes as plain text strings, some containing multiple dots).

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.232.62'
$ws.Range('E2').Value = '  +1.31%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.817.78'
$ws.Range('E3').Value = '  -1.92%  '

# Row 4
$ws.Range('E4').Value = '  +0.31%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.94'
$ws.Range('E5').Value = '  +2.31%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.614'
$ws.Range('E6').Value = '  +0.55%  '

# Row 7
$ws.Range('E7').Value = '  +0.30%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.03'
$ws.Range('E8').Value = '  +0.27%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.327'
$ws.Range('E9').Value = '  +7.45%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0686'
$ws.Range('E10').Value = '  -0.75%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0999'
$ws.Range('E11').Value = '  -0.33%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.081.99'
$ws.Range('E12').Value = '  -1.81%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.820.21'
$ws.Range('E13').Value = '  -1.72%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.16'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.661'
$ws.Range('E15').Value = '  +0.19%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.67'
$ws.Range('E16').Value = '  -1.54%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.165.40'
$ws.Range('E17').Value = '  +1.27%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.73'
$ws.Range('E18').Value = '  +1.02%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0794'
$ws.Range('E19').Value = '  +0.65%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.09'
$ws.Range('E20').Value = '  -1.98%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.99'
$ws.Range('E21').Value = '  -1.41%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.66'
$ws.Range('E22').Value = '  -2.31%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.17%  '

# Row 24
$ws.Range('E24').Value = '  +3.81%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.31'
$ws.Range('E25').Value = '  -0.14%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.86'
$ws.Range('E26').Value = '  -0.98%  '

# Row 27
$ws.Range('E27').Value = '  -1.70%  '

# Row 28
$ws.Range('E28').Value = '  -1.04%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.61'
$ws.Range('E29').Value = '  +20.39%  '

# Row 30
$ws.Range('E30').Value = '  +0.35%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.333.73'
$ws.Range('E31').Value = '  +37.21%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.04'
$ws.Range('E32').Value = '  +2.52%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0556'
$ws.Range('E33').Value = '  +4.22%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.00'
$ws.Range('E34').Value = '  -0.18%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.78'
$ws.Range('E35').Value = '  -6.43%  '

# Row 36
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.14'
$ws.Range('E36').Value = '  +6.32%  '

# Row 37
$ws.Range('B37').Value = 'Aave'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '93.27'
$ws.Range('E37').Value = '  +1.96%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.684'
$ws.Range('E38').Value = '  +2.05%  '

# Row 39
$ws.Range('E39').Value = '  +0.47%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.28'
$ws.Range('E40').Value = '  +2.72%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.312.20'
$ws.Range('E41').Value = '  -2.46%  '

# Row 42
$ws.Range('E42').Value = '  -3.13%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '14.73'
$ws.Range('E43').Value = '  -0.58%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.31'
$ws.Range('E44').Value = '  -5.33%  '

# Row 46
$ws.Range('E46').Value = '  -2.69%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.34'
$ws.Range('E47').Value = '  +5.04%  '

# Row 48
$ws.Range('E48').Value = '  -1.01%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.998.44'
$ws.Range('E49').Value = '  -0.90%  '

# Row 50
$ws.Range('E50').Value = '  +0.19%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0648'
$ws.Range('E51').Value = '  +5.19%  '
